$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "aliii"
$ws.Range("B2").Value = "ali12@ee.com"
$ws.Range("D2").Value = "pass22"
$ws.Range("E2").Value = "pass22"

# Row 3
$ws.Range("A3").Value = "stt1"
$ws.Range("B3").Value = "stt1@ee.com"
$ws.Range("D3").Value = "pass22"
$ws.Range("E3").Value = "pass22"

# Row 4
$ws.Range("A4").Value = "stt11"
$ws.Range("B4").Value = "stt11@o.com"
$ws.Range("D4").Value = "pass22"
$ws.Range("E4").Value = "pass22"

# Row 5
$ws.Range("A5").Value = "stt111"
$ws.Range("B5").Value = "stt111@e.com"
$ws.Range("D5").Value = "pass22"
$ws.Range("E5").Value = "pass22"

# Update the active selection to E4 as recorded in the saved view state
$ws.Range("E4").Select()

$wb.Save()
